$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 82 (Aysén block), shifting the
# existing rows 82:87 down to 83:88.
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row with the new candidate record.
$ws.Range("A82").Value = 87
$ws.Range("B82").Value = 13
$ws.Range("C82").Value = 11
$ws.Range("D82").Value = "Aysén"
$ws.Range("E82").Value = "Jorge Sepúlveda"
$ws.Range("F82").Value = "Independiente"
$ws.Range("G82").Value = "IND"
$ws.Range("H82").Value = 0

# Re-apply the existing sort (same keys, grown by one row) so the sheet's
# remembered sort range/conditions grow to cover the newly inserted row.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C89"), 0, 1, 0, 0) | Out-Null
$ws.Sort.SortFields.Add($ws.Range("F2:F89"), 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:H89"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Match the author's final selection.
$ws.Range("E4").Select()
